# Fix author-name typo throughout the deck: 吳建宏 -> 吳建鋐
#
# Build the old/new name from Unicode code points (avoids any
# source-encoding pitfalls) using string interpolation - NOT '+', since
# this engine treats [char] + [char] as numeric addition rather than
# concatenation.
$wu    = [char]21555   # 吳
$jian  = [char]24314   # 建
$hong  = [char]23439   # 宏  (old / typo)
$hong2 = [char]37584   # 鋐  (new / corrected)

$oldName = "$wu$jian$hong"
$newName = "$wu$jian$hong2"
$oldLen  = $oldName.Length

function Fix-ShapeText($shape) {
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $searchStart = 0
    while ($true) {
        $idx = $full.IndexOf($oldName, $searchStart)
        if ($idx -lt 0) { break }
        # TextRange.Characters(Start, Length) is 1-based and lets us
        # replace just the matched substring in place, leaving
        # surrounding runs/formatting untouched.
        $sub = $tr.Characters($idx + 1, $oldLen)
        $sub.Text = $newName
        $searchStart = $idx + $oldLen
    }
}

$p = $ppt.ActivePresentation

# Slide 1: title slide author list ("吳建宏" on its own line)
$s1 = $p.Slides.Item(1)
Fix-ShapeText $s1.Shapes.Item(1)

# Slide 8: contribution breakdown ("發想:吳建宏、林佛兒", "GUI程式:吳建宏",
# "測試與調整:吳建宏、林佛兒")
$s8 = $p.Slides.Item(8)
Fix-ShapeText $s8.Shapes.Item(2)
